$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 128
$ws.Range("B1").Value = 128
$ws.Range("C1").Value = 128
$ws.Range("D1").Value = 53.585013452169022

$ws.Range("A2").Value = 92
$ws.Range("C2").Value = 80
$ws.Range("D2").Value = 19.453894549302518
$ws.Range("E2").Value = 44.726077651902841
$ws.Range("F2").Value = -22.042905540037182

$ws.Range("B3").Value = 84
$ws.Range("C3").Value = 241
$ws.Range("D3").Value = 42.290701884164996
$ws.Range("E3").Value = 41.31476307866788
$ws.Range("F3").Value = -83.398685764479751

$ws.Range("B4").Value = 154
$ws.Range("C4").Value = 48
$ws.Range("D4").Value = 55.403475712623731
$ws.Range("E4").Value = -56.328173879348576
$ws.Range("F4").Value = 44.306292693830699

$ws.Range("B5").Value = 251
$ws.Range("D5").Value = 96.103284087075764
$ws.Range("E5").Value = -19.650215533887405
$ws.Range("F5").Value = 93.721096841030899

$ws.Range("A6").Value = 239
$ws.Range("B6").Value = 144
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 68.242750803252676
$ws.Range("E6").Value = 28.232703032031846
$ws.Range("F6").Value = 73.659648978326516

$ws.Range("B8").Value = 255
$ws.Range("C8").Value = 255
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
